$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Insert a "_GoBack" bookmark at the very start of the document (before
#    the "Título" run of the first paragraph). An empty Range placed exactly
#    at document position 0 behaves oddly for Bookmarks.Add (the resulting
#    bookmarkEnd ends up on the following paragraph), so we temporarily
#    insert a one-character placeholder at the start, wrap *that* character
#    with the bookmark, then delete the placeholder character again. This
#    leaves bookmarkStart/bookmarkEnd adjacent to each other, right before
#    the first run, exactly as in the target document.
# ---------------------------------------------------------------------------
$startR = $d.Range(0, 0)
$startR.InsertBefore("X")
$placeholder = $d.Range(0, 1)
$d.Bookmarks.Add("_GoBack", $placeholder)
$d.Range(0, 1).Delete()

# ---------------------------------------------------------------------------
# 2. The body text has "Etia" + a stray "_GoBack" bookmark + "m" split across
#    three runs (an artifact of earlier editing). Collapse that back into a
#    single "Etiam" run and drop the bookmark. Setting Range.Text is a no-op
#    when the replacement text is identical to the existing text, so first
#    swap in a scratch value and then swap in the real word; this forces the
#    underlying run (and the markup nested inside the old range, including
#    the bookmark) to be rebuilt as one plain run.
# ---------------------------------------------------------------------------
$etia = $d.Content
[void]$etia.Find.Execute("Etia", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
[void]$etia.MoveEnd(1, 1)
$etia.Text = "%%TMP%%"

$etiam = $d.Content
[void]$etiam.Find.Execute("%%TMP%%", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$etiam.Text = "Etiam"

# ---------------------------------------------------------------------------
# 3. Reformat the three TOC paragraph styles (Sumário 1/2/3): space-after
#    goes from 3pt (60 twips) to 6pt (120 twips).
# ---------------------------------------------------------------------------
foreach ($styleName in @("Sumrio1", "Sumrio2", "Sumrio3")) {
    $style = $d.Styles($styleName)
    $style.ParagraphFormat.SpaceAfter = 6
}
